$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- D4: "An ERD" was retyped, dropping the leading non-breaking space so
# only a single plain space remains before the text. Keep the original
# two-run rich formatting (small Times New Roman lead-in space, then the
# Calibri label text).
$ws.Range("D4").Value = " An ERD"
$leadSpace = $ws.Range("D4").Characters(1, 1)
$leadSpace.Font.Size = 7
$leadSpace.Font.Name = "Times New Roman"
$label = $ws.Range("D4").Characters(2, 6)
$label.Font.Size = 11
$label.Font.Name = "Calibri"

# --- Actual Hours (column F) updates ---
# SQL database script row: actual hours revised from 15 to 20.
$ws.Range("F5").Value = 20

# These rows previously had no Actual Hours recorded; fill them in.
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 1

# Total row for Actual Hours (entered as a plain value, not a formula).
$ws.Range("F13").Value = 34.5

# Update the active selection left over from the editing session.
[void]$ws.Range("G11").Select()
